{"js": "// Office.js (Word JavaScript API) script\n// Implements: \"Oskar Agrees to the WBA\"\n//  1) In the WBA table, row \"7: Reset\", append \" & Shyam\" to the two cells\n//     (Implementation, Testing & Reviewing) that currently just say \"Oskar\".\n//  2) After the final \"I agree to this WBA \u2013 Shyam\" paragraph, add a new\n//     paragraph \"I agree to this WBA - Oskar\" with matching formatting.\n\nconst body = context.document.body;\n\n// --- 1) Update the \"7: Reset\" row in the (only) table ----------------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nrows.items.forEach((r) => r.cells.load(\"items\"));\nawait context.sync();\n\n// Locate the row whose first cell starts with \"7: Reset\".\nconst firstCellBodies = rows.items.map((row) => row.cells.items[0].body);\nfirstCellBodies.forEach((b) => b.load(\"text\"));\nawait context.sync();\n\nlet targetRow = null;\nfor (let i = 0; i < rows.items.length; i++) {\n  const firstCellText = firstCellBodies[i].text;\n  if (firstCellText && firstCellText.trim().startsWith(\"7: Reset\")) {\n    targetRow = rows.items[i];\n    break;\n  }\n}\n\nif (!targetRow) {\n  throw new Error(\"Could not find the '7: Reset' row in the table.\");\n}\n\n// Columns 1 (Implementation) and 2 (Testing & Reviewing) need \" & Shyam\"\n// appended right after the existing \"Oskar\" text. Inserting at the end of\n// the cell's last paragraph (without touching font properties afterwards)\n// lets the new text inherit the formatting of the run it is appended to.\nconst colIndexes = [1, 2];\nfor (const idx of colIndexes) {\n  const cell = targetRow.cells.items[idx];\n  const paras = cell.body.paragraphs;\n  paras.load(\"items\");\n  await context.sync();\n\n  const lastPara = paras.items[paras.items.length - 1];\n  lastPara.insertText(\" & Shyam\", \"End\");\n}\n\nawait context.sync();\n\n// --- 2) Add a new \"I agree...\" paragraph for Oskar -------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastBodyPara = paragraphs.items[paragraphs.items.length - 1];\nlastBodyPara.insertParagraph(\"I agree to this WBA - Oskar\", \"After\");\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script\n# Implements: \"Oskar Agrees to the WBA\"\n#  1) In the WBA table, row \"7: Reset\", append \" & Shyam\" to the two cells\n#     (Implementation, Testing & Reviewing) that currently just say \"Oskar\".\n#  2) After the final \"I agree to this WBA \u2013 Shyam\" paragraph, add a new\n#     paragraph \"I agree to this WBA - Oskar\" with the same formatting.\n\n$d = $word.ActiveDocument\n\n# --- 1) Update the \"7: Reset\" row in the (only) table ---------------------\n$tbl = $d.Tables.Item(1)\n\n$targetRow = 0\nfor ($i = 1; $i -le $tbl.Rows.Count; $i++) {\n    $rowLabel = $tbl.Cell($i, 1).Range.Text\n    if ($rowLabel.StartsWith(\"7: Reset\")) {\n        $targetRow = $i\n        break\n    }\n}\n\nif ($targetRow -eq 0) {\n    throw \"Could not find the '7: Reset' row in the table.\"\n}\n\nforeach ($col in 2, 3) {\n    $cell = $tbl.Cell($targetRow, $col)\n    $cellRng = $cell.Range\n    # Cell.Range includes the trailing paragraph mark; back off one\n    # character so the new text lands right after \"Oskar\" and inherits\n    # its run formatting (Arial, 12pt).\n    $cellRng.MoveEnd(1, -1)\n    $cellRng.InsertAfter(\" & Shyam\")\n}\n\n# --- 2) Add a new \"I agree...\" paragraph for Oskar -------------------------\n# Use $d.Content.Paragraphs (rather than $d.Paragraphs) to get a fresh,\n# document-order collection that is unaffected by the table access above.\n$paragraphs = $d.Content.Paragraphs\n$lastPara = $paragraphs.Item($paragraphs.Count)\n$lastRange = $lastPara.Range\n$lastRange.InsertParagraphAfter()\n\n$paragraphs = $d.Content.Paragraphs\n$newPara = $paragraphs.Item($paragraphs.Count)\n$newPara.Range.InsertAfter(\"I agree to this WBA - Oskar\")\n"}
